# Generate Report for handoff
# Update "Latest Handoff Datetime" (column D) for the 4th data row (row 5,
# corresponding to the a409dc41-8675-45ae-bb8a-3121e5df4d6a entry) on both
# the "zh-cn" and "de-de" sheets with the new handoff timestamps produced
# by the latest report generation run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-13 11:24:55"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-13 11:25:20"
